# Update column G ("K") values for rows 2-16 on Sheet1.
# The data file was regenerated to use "K" (strikeouts?) values instead of
# the previous "Strike#" values, so only the G column numbers change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 2
    7  = 0
    8  = 3
    9  = 2
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}

$wb.Save()
